$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = '30.604.79'
$ws.Range("E2").Value = '  +0.71%  '

$ws.Range("D3").Value = '2.116.30'
$ws.Range("E3").Value = '  +0.24%  '

$ws.Range("E4").Value = '  +0.64%  '

$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '350.55'
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = '  +4.88%  '

$ws.Range("E6").Value = '  +0.70%  '

$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = '0.5258'
$ws.Range("D7").Style = "Normal"
$ws.Range("E7").Value = '  +0.43%  '

$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = '0.4510'
$ws.Range("D8").Style = "Normal"
$ws.Range("E8").Value = '  -1.13%  '

$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = '54.60'
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = '  +2.04%  '

$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = '0.09048'
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = '  +1.42%  '

$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = '1.178'
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = '  -0.17%  '

$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = '24.55'
$ws.Range("D12").Style = "Normal"
$ws.Range("E12").Value = '  +0.45%  '

$ws.Range("D13").Value = '2.100.46'
$ws.Range("E13").Value = '  +0.01%  '

$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = '6.835'
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = '  -0.08%  '

$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = '8.068'
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = '  +0.27%  '

$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = '102.61'
$ws.Range("D16").Style = "Normal"
$ws.Range("E16").Value = '  +6.16%  '

$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = '0.00001173'
$ws.Range("D17").Style = "Normal"
$ws.Range("E17").Value = '  +3.26%  '

$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = '1.011'
$ws.Range("D18").Style = "Normal"
$ws.Range("E18").Value = '  +0.66%  '

$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = '0.06714'
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = '  +1.06%  '

$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = '19.45'
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = '  +0.78%  '

$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = '1.009'
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = '  +0.73%  '

$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = '6.309'
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = '  -1.03%  '

$ws.Range("D23").Value = '30.661.27'
$ws.Range("E23").Value = '  +0.61%  '

$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = '12.83'
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = '  +3.17%  '

$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = '2.382'
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = '  +0.70%  '

$ws.Range("D26").Value = '2.357.44'
$ws.Range("E26").Value = '  +0.67%  '

$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = '22.50'
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = '  +0.41%  '

$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = '165.22'
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = '  +0.82%  '

$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = '2.552'
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = '  -0.76%  '

$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = '137.15'
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = '  +3.19%  '

$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = '1.195'
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = '  -3.89%  '

$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = '0.1078'
$ws.Range("D32").Style = "Normal"

$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = '1.667'
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = '  -2.53%  '

$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = '6.380'
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = '  -0.04%  '

$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = '4.019'
$ws.Range("D35").Style = "Normal"

$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = '10.41'
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = '  -1.36%  '

$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = '5.920'
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = '  +5.45%  '

$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = '0.02648'
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = '  +2.28%  '

$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = '0.06868'
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = '  +0.48%  '

$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = '0.2316'
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = '  +0.26%  '

$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = '12.57'
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = '  -1.62%  '

$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = '0.6897'
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = '  -0.49%  '

$ws.Range("E43").Value = '  +1.91%  '

$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = '14.73'
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = '  +4.40%  '

$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = '2.334'
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = '  -0.96%  '

$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = '0.6464'
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = '  +0.99%  '

$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = '3.755'
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = '  +2.65%  '

$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = '0.00000000363'
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = '  +2.92%  '

$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = '1.254'
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = '  +0.33%  '

$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = '0.07295'
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = '  +2.12%  '

$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = '82.61'
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = '  -1.07%  '
